$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Row 2  <=  original row 11 data
$ws.Range("A2").Value = 89596126
$ws.Range("B2").Value = 56411
$ws.Range("D2").Value = "NT"
$ws.Range("E2").Value = 100049
$ws.Range("F2").Value = "Spillkråka"
$ws.Range("G2").Value = "Dryocopus martius"
$ws.Range("H2").Value = "(Linnaeus, 1758)"
$ws.Range("P2").Value = "Strömsfjällvallen, Hjd"
$ws.Range("Q2").Value = 444929.0050177791
$ws.Range("R2").Value = 6928327.074997591
$ws.Range("S2").Value = 10
$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value = "2020-09-25"
$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value = "2020-09-25"
$ws.Range("AC2").Value = "Födosökshack"
$ws.Range("AW2").Value = "Erland Lindblad"
$ws.Range("AX2").Value = "Jan Henriksson"
$ws.Range("AY2").Value = "Kontinuitetsskogar och skogar med höga naturvärden ovan och i nära anslutning till fjällnära gränsen"

# Row 3  <=  original row 5 data
$ws.Range("A3").Value = 89596129
$ws.Range("B3").Value = 76909
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 6437
$ws.Range("F3").Value = "Blanksvart spiklav"
$ws.Range("G3").Value = "Calicium denigratum"
$ws.Range("H3").Value = "(Vain.) Tibell"
$ws.Range("P3").Value = "Strömsfjällvallen, Hjd"
$ws.Range("Q3").Value = 445032.0268228107
$ws.Range("R3").Value = 6928535.7939387
$ws.Range("S3").Value = 10
$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value = "2020-09-25"
$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value = "2020-09-25"
$ws.Range("AC3").ClearContents()
$ws.Range("AW3").Value = "Erland Lindblad"
$ws.Range("AX3").Value = "Jan Henriksson"
$ws.Range("AY3").Value = "Kontinuitetsskogar och skogar med höga naturvärden ovan och i nära anslutning till fjällnära gränsen"

# Row 4  <=  original row 12 data
$ws.Range("A4").Value = 89596127
$ws.Range("B4").Value = 73693
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 6440
$ws.Range("F4").Value = "Vitgrynig nållav"
$ws.Range("G4").Value = "Chaenotheca subroscida"
$ws.Range("H4").Value = "(Eitner) Zahlbr."
$ws.Range("P4").Value = "Strömsfjällvallen, Hjd"
$ws.Range("Q4").Value = 444967.7734563763
$ws.Range("R4").Value = 6928430.952647353
$ws.Range("S4").Value = 10
$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value = "2020-09-25"
$ws.Range("AA4").NumberFormat = "@"
$ws.Range("AA4").Value = "2020-09-25"
$ws.Range("AC4").ClearContents()
$ws.Range("AW4").Value = "Erland Lindblad"
$ws.Range("AX4").Value = "Jan Henriksson"
$ws.Range("AY4").Value = "Kontinuitetsskogar och skogar med höga naturvärden ovan och i nära anslutning till fjällnära gränsen"

# Row 5  <=  original row 6 data
$ws.Range("A5").Value = 89596128
$ws.Range("B5").Value = 77506
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 6425
$ws.Range("F5").Value = "Garnlav"
$ws.Range("G5").Value = "Alectoria sarmentosa"
$ws.Range("H5").Value = "(Ach.) Ach."
$ws.Range("P5").Value = "Strömsfjällvallen, Hjd"
$ws.Range("Q5").Value = 445020.0425569176
$ws.Range("R5").Value = 6928539.228960108
$ws.Range("S5").Value = 10
$ws.Range("Y5").NumberFormat = "@"
$ws.Range("Y5").Value = "2020-09-25"
$ws.Range("AA5").NumberFormat = "@"
$ws.Range("AA5").Value = "2020-09-25"
$ws.Range("AC5").ClearContents()
$ws.Range("AW5").Value = "Erland Lindblad"
$ws.Range("AX5").Value = "Jan Henriksson"
$ws.Range("AY5").Value = "Kontinuitetsskogar och skogar med höga naturvärden ovan och i nära anslutning till fjällnära gränsen"

# Row 6  <=  original row 7 data
$ws.Range("A6").Value = 94995564
$ws.Range("B6").Value = 95525
$ws.Range("D6").Value = "LC"
$ws.Range("E6").Value = 221941
$ws.Range("F6").Value = "Plattlummer"
$ws.Range("G6").Value = "Lycopodium complanatum"
$ws.Range("H6").Value = "L."
$ws.Range("P6").Value = "Härjedalen, Hjd"
$ws.Range("Q6").Value = 445086.0189850244
$ws.Range("R6").Value = 6928496.057011075
$ws.Range("S6").Value = 25
$ws.Range("Y6").NumberFormat = "@"
$ws.Range("Y6").Value = "2021-07-20"
$ws.Range("AA6").NumberFormat = "@"
$ws.Range("AA6").Value = "2021-07-20"
$ws.Range("AC6").Value = "Tallskog, ås."
$ws.Range("AW6").Value = "Jens Johannesson"
$ws.Range("AX6").Value = "Jens Johannesson"
$ws.Range("AY6").ClearContents()

# Row 7  <=  original row 3 data
$ws.Range("A7").Value = 89596109
$ws.Range("B7").Value = 78570
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 2081
$ws.Range("F7").Value = "Skrovellav"
$ws.Range("G7").Value = "Lobaria scrobiculata"
$ws.Range("H7").Value = "(Scop.) DC."
$ws.Range("P7").Value = "Strömsfjällvallen, Hjd"
$ws.Range("Q7").Value = 445260.1075701897
$ws.Range("R7").Value = 6928606.027293501
$ws.Range("S7").Value = 10
$ws.Range("Y7").NumberFormat = "@"
$ws.Range("Y7").Value = "2020-09-25"
$ws.Range("AA7").NumberFormat = "@"
$ws.Range("AA7").Value = "2020-09-25"
$ws.Range("AC7").ClearContents()
$ws.Range("AW7").Value = "Erland Lindblad"
$ws.Range("AX7").Value = "Jan Henriksson"
$ws.Range("AY7").Value = "Kontinuitetsskogar och skogar med höga naturvärden ovan och i nära anslutning till fjällnära gränsen"

# Row 11  <=  original row 2 data
$ws.Range("A11").Value = 94996219
$ws.Range("B11").Value = 78570
$ws.Range("D11").Value = "NT"
$ws.Range("E11").Value = 2081
$ws.Range("F11").Value = "Skrovellav"
$ws.Range("G11").Value = "Lobaria scrobiculata"
$ws.Range("H11").Value = "(Scop.) DC."
$ws.Range("P11").Value = "Härjedalen, Hjd"
$ws.Range("Q11").Value = 445370.7671139772
$ws.Range("R11").Value = 6928604.672176878
$ws.Range("S11").Value = 25
$ws.Range("Y11").NumberFormat = "@"
$ws.Range("Y11").Value = "2021-07-20"
$ws.Range("AA11").NumberFormat = "@"
$ws.Range("AA11").Value = "2021-07-20"
$ws.Range("AC11").Value = "På jättesälg ca 45 cm diameter."
$ws.Range("AW11").Value = "Jens Johannesson"
$ws.Range("AX11").Value = "Jens Johannesson"
$ws.Range("AY11").ClearContents()

# Row 12  <=  original row 4 data
$ws.Range("A12").Value = 94996005
$ws.Range("B12").Value = 78570
$ws.Range("D12").Value = "NT"
$ws.Range("E12").Value = 2081
$ws.Range("F12").Value = "Skrovellav"
$ws.Range("G12").Value = "Lobaria scrobiculata"
$ws.Range("H12").Value = "(Scop.) DC."
$ws.Range("P12").Value = "Härjedalen, Hjd"
$ws.Range("Q12").Value = 445261.8150698114
$ws.Range("R12").Value = 6928597.212872105
$ws.Range("S12").Value = 25
$ws.Range("Y12").NumberFormat = "@"
$ws.Range("Y12").Value = "2021-07-20"
$ws.Range("AA12").NumberFormat = "@"
$ws.Range("AA12").Value = "2021-07-20"
$ws.Range("AC12").Value = "Mkt gammal sälg, fin skog på åsar."
$ws.Range("AW12").Value = "Jens Johannesson"
$ws.Range("AX12").Value = "Jens Johannesson"
$ws.Range("AY12").ClearContents()
